$wb = $excel.ActiveWorkbook

# "Raw Data" sheet: add row 70 with the new product's raw seconds value
$wsRaw = $wb.Worksheets.Item("Raw Data")
$wsRaw.Range("A70").Value = "110256_D"
$wsRaw.Range("B70").Value = 1404.184306

# "Results" sheet: add row 49 with the aggregated product code and formatted time
# (leading apostrophe forces the numeric-looking code to stay text, matching
#  the other product-code cells in this column)
$wsResults = $wb.Worksheets.Item("Results")
$wsResults.Range("A49").Value = "'110256"
$wsResults.Range("B49").Value = "0:23:24.184306000000106"
